$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fire-station")
$ws.Select()
Write-Output "hello"
